# The document's TOC heading bookmarks (Word's internal "_Toc..." hidden
# bookmarks) were regenerated with new ids by Word when the document was
# re-saved ("prima di rimuovere footer"). Bookmark.Name is read-only in the
# Word object model, so each bookmark is recreated (same range, new name)
# rather than renamed in place.

$d = $word.ActiveDocument

# Map of old hidden "_Toc" bookmark name -> new name, in document order.
$renames = [ordered]@{
    "_Toc5791590" = "_Toc5794235"
    "_Toc5791591" = "_Toc5794236"
    "_Toc5791592" = "_Toc5794237"
    "_Toc5791593" = "_Toc5794238"
    "_Toc5791594" = "_Toc5794239"
    "_Toc5791595" = "_Toc5794240"
}

foreach ($oldName in $renames.Keys) {
    $newName = $renames[$oldName]

    if ($d.Bookmarks.Exists($oldName)) {
        $bm = $d.Bookmarks($oldName)
        $rng = $bm.Range
        $bm.Delete()
        $d.Bookmarks.Add($newName, $rng)
    }
}
